$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 and 35: coin name/link swap (OKB <-> InjectiveProtocol), plus new price/volume values
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.45"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.28"
$ws.Range("E35").Value = "  +0.67%  "

# Price (D) and Volume(1h) (E) updates for remaining rows
# D-column cells are forced to Text format ("@") since the original
# values are stored as text (inline strings), and several new values
# (e.g. "1.00", "380.26") would otherwise be auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.073.50"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.960.45"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.26"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.25"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.545"
$ws.Range("E7").Value = "  +1.74%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.588"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.63"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0852"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.424.42"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.39"
$ws.Range("E14").Value = "  +2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.76"
$ws.Range("E15").Value = "  +5.57%  "
$ws.Range("E16").Value = "  +70.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.959.30"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.00"
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.161.11"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.09"
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.33"
$ws.Range("E23").Value = "  +16.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.20"
$ws.Range("E24").Value = "  +2.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.77"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.94"
$ws.Range("E26").Value = "  -2.74%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.89"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.04"
$ws.Range("E30").Value = "  -10.77%  "
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("E32").Value = "  +6.20%  "
$ws.Range("E33").Value = "  +6.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0436"
$ws.Range("E36").Value = "  -4.18%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  +10.89%  "
$ws.Range("E39").Value = "  +1.96%  "
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.52"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.76"
$ws.Range("E43").Value = "  +3.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.73"
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("E45").Value = "  +10.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.272"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.064.47"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0320"
$ws.Range("E50").Value = "  -9.11%  "
$ws.Range("E51").Value = "  +6.99%  "
